$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headline = 'mueller industry inc nysemli share sold by convergence investment partner llc convergence investment partner llc sell share of commvault system inc nasdaqcvlt kimberly s greene sell share of southern company the nyseso stock kellanova k to release earnings on thursday oak thistle llc ha million stake in booking holding inc nasdaqbkng venturi wealth management llc ha million stake in ishares core sp midcap etf nysearcaijh share in quidelortho corporation nasdaqqdel purchased by sym financial corp ishares msci eafe growth etf batsefg stake lowered by venturi wealth management llc venturi wealth management llc increase stock position in spotify technology nysespot venturi wealth management llc buy share of cullenfrost banker inc nysecfr arista network inc nyseanet share sold by venturi wealth management llc y intercept hong kong ltd lower stock holding in the traveler company inc nysetrv y intercept hong kong ltd invests million in moody corporation nysemco venturi wealth management llc buy share of valmont industry inc nysevmi globe life inc nysegl position lifted by venturi wealth management llc venturi wealth management llc sell share of ishares mb etf nasdaqmbb venturi wealth management llc purchase share of brookfield corporation nysebn sym financial corp make new investment in cv health corporation nysecvs share in ugi corporation nyseugi purchased by sym financial corp venturi wealth management llc purchase share of armstrong world industry inc nyseawi universal beteiligungs und servicegesellschaft mbh purchase share of biogen inc nasdaqbiib universal beteiligungs und servicegesellschaft mbh sell share of cognizant technology solution corporation nasdaqctsh venturi wealth management llc ha million stock position in reliance inc nysers universal beteiligungs und servicegesellschaft mbh ha million stock position in biogen inc nasdaqbiib venturi wealth management llc raise holding in lharris technology inc nyselhx venturi wealth management llc trim stake in ishares u technology etf nysearcaiyw sym financial corp purchase share of ugi corporation nyseugi venturi wealth management llc raise stock position in valmont industry inc nysevmi y intercept hong kong ltd sell share of synopsys inc nasdaqsnps intel corporation nasdaqintc share sold by sym financial corp brookfield corporation nysebn stock holding lifted by venturi wealth management llc y intercept hong kong ltd purchase share of vanguard mega cap growth etf nysearcamgk sym financial corp ha stock holding in philip morris international inc nysepm lincoln national corporation nyselnc share acquired by sym financial corp venturi wealth management llc ha position in armstrong world industry inc nyseawi sym financial corp ha stock position in steel dynamic inc nasdaqstld venturi wealth management llc lower stake in morgan stanley nysems cognizant technology solution corporation nasdaqctsh share sold by universal beteiligungs und servicegesellschaft mbh cullenfrost banker inc nysecfr share acquired by venturi wealth management llc mcdonalds to test cosmcsinspired drink at more than restaurant pfg advisor ha million stake in altshares merger arbitrage etf nysearcaarb choreo llc reduces stock position in first financial bancorp nasdaqffbc pfg advisor buy share of ishares year treasury bond etf nasdaqshy choreo llc sell share of brown brown inc nysebro t rowe price u high yield etf nysearcathyf stake raised by pfg advisor american international group inc nyseaig holding lowered by choreo llc choreo llc decrease stake in ishares year tip bond etf nysearcastip danaher corporation nysedhr share sold by choreo llc y intercept hong kong ltd acquires share of tetra tech inc nasdaqttek schwab u midcap etf nysearcaschm stock position lowered by choreo llc choreo llc ha million stock position in spdr dow jones industrial average etf trust nysearcadia the pnc financial service group inc nysepnc share sold by choreo llc choreo llc cut stake in autozone inc nyseazo equinix inc nasdaqeqix share purchased by pfg advisor y intercept hong kong ltd take million position in tyson food inc nysetsn choreo llc sell share of dr horton inc nysedhi y intercept hong kong ltd cut stake in nrg energy inc nysenrg ishares esg aware msci eafe etf nasdaqesgd share sold by choreo llc choreo llc ha million position in darden restaurant inc nysedri highpoint advisor group llc lower stake in tetra tech inc nasdaqttek choreo llc ha million stock holding in vanguard dividend appreciation etf nysearcavig choreo llc sell share of parkerhannifin corporation nyseph y intercept hong kong ltd acquires share of axon enterprise inc nasdaqaxon choreo llc sell share of occidental petroleum corporation nyseoxy choreo llc grows stake in micron technology inc nasdaqmu y intercept hong kong ltd boost stake in eastman chemical company nyseemn share in eagle bancorp inc nasdaqegbn purchased by readystate asset management lp readystate asset management lp take position in cnx resource corporation nysecnx readystate asset management lp purchase share of cnx resource corporation nysecnx y intercept hong kong ltd sell share of illumina inc nasdaqilmn best ai video generation tool keysight technology q earnings what to expect nvidia address ai chip smuggling say bootleg datacenters are a losing proposition want to sound more assertive ditch this word phrase say trial lawyer and communication expert readystate asset management lp take position in western asset high yield defined opportunity fund inc nysehyi ishares tip bond etf nysearcatip share purchased by proficio capital partner llc readystate asset management lp increase stock holding in american axle manufacturing holding inc nyseaxl readystate asset management lp ha position in nov inc nysenov proficio capital partner llc decrease holding in novartis ag nysenvs proficio capital partner llc sell share of kimberlyclark co nysekmb share in franklin limited duration income trust nyseamericanftf purchased by readystate asset management lp klp kapitalforvaltning a increase stock position in old dominion freight line inc nasdaqodfl proficio capital partner llc ha stock holding in otis worldwide corporation nyseotis readystate asset management lp grows position in halliburton company nysehal united parcel service inc nyseups share bought by proficio capital partner llc proficio capital partner llc sell share of unum group nyseunm share in rogers communication inc nyserci bought by proficio capital partner llc proficio capital partner llc ha stake in under armour inc nyseuaa readystate asset management lp acquires share of halliburton company nysehal readystate asset management lp take position in columbia banking system inc nasdaqcolb proficio capital partner llc increase stake in united parcel service inc nyseups proficio capital partner llc purchase share of rogers communication inc nyserci proficio capital partner llc sell share of exlservice holding inc nasdaqexls proficio capital partner llc ha position in mondelez international inc nasdaqmdlz medpace holding inc nasdaqmedp share purchased by proficio capital partner llc readystate asset management lp take position in virtu total return fund inc nyseztr proficio capital partner llc sell share of mondelez international inc nasdaqmdlz proficio capital partner llc purchase share of norfolk southern corporation nysensc proficio capital partner llc sell share of nisource inc nyseni proficio capital partner llc ha stake in norfolk southern corporation nysensc which fabric name is used for the tissue that grows on deer antler k street rake in record cash thanks to trump rnc chair michael whatley to run for senate in north carolina with trump support politico lowell city council move toward a green energy future venturi wealth management llc reduces stock position in ishares u technology etf nysearcaiyw sym financial corp purchase new share in draftkings inc nasdaqdkng venturi wealth management llc grows position in amdocs limited nasdaqdox venturi wealth management llc ha holding in lyft inc nasdaqlyft sym financial corp ha stake in steel dynamic inc nasdaqstld venturi wealth management llc ha stock holding in lharris technology inc nyselhx sym financial corp make new investment in axis capital holding limited nyseaxs sym financial corp reduces stake in intel corporation nasdaqintc sym financial corp acquires share of lincoln national corporation nyselnc y intercept hong kong ltd increase position in raymond james financial inc nyserjf universal beteiligungs und servicegesellschaft mbh cut stock holding in dell technology inc nysedell venturi wealth management llc sell share of morgan stanley nysems philip morris international inc nysepm share sold by sym financial corp universal beteiligungs und servicegesellschaft mbh ha million stock holding in royal caribbean cruise ltd nysercl y intercept hong kong ltd sell share of synopsys inc nasdaqsnps readystate asset management lp invests in fulton financial corporation nasdaqfult y intercept hong kong ltd ha million stock holding in illumina inc nasdaqilmn trump set to head to fed keeping heat on chair jerome powell letter democracy alive and well under trump readystate asset management lp invests in columbia banking system inc nasdaqcolb proficio capital partner llc sell share of genuine part company nysegpc klp kapitalforvaltning a cut position in electronic art inc nasdaqea readystate asset management lp purchase new holding in virtu total return fund inc nyseztr proficio capital partner llc ha stock position in seritage growth property nysesrg proficio capital partner llc purchase share of medpace holding inc nasdaqmedp exlservice holding inc nasdaqexls share sold by proficio capital partner llc proficio capital partner llc ha holding in crowdstrike nasdaqcrwd klp kapitalforvaltning a ha million stock position in lennar corporation nyselen jeffrey p kenefick purchase share of chemung financial corp nasdaqchmg stock klp kapitalforvaltning a acquires share of coinbase global inc nasdaqcoin insider buying ioneer ltd asxinr insider purchase a in stock klp kapitalforvaltning a boost holding in sea limited sponsored adr nysese klp kapitalforvaltning a increase stock position in t rowe price group inc nasdaqtrow nxt energy solution inc tsesfd insider sell c in stock insider buying innovator u equity ultra buffer etf november batsunov ceo buy in stock moishe gubin buy share of strawberry field reit inc nasdaqstrw stock netapp inc nasdaqntap share sold by klp kapitalforvaltning a letter duty to care for the poor contrasting paysign nasdaqpays and evertec nyseevtc kinross gold fy eps forecast boosted by stifel canada viking nysevik price target raised to ross store nasdaqrost stock unloaded rep jefferson shreve brokerage set prologis inc nysepld price target at ge aerospace nysege stock unloaded rep jefferson shreve akamai technology nasdaqakam stock unloaded rep jefferson shreve brokerage set juniper network inc nysejnpr price target at q eps estimate for formation mnrls lowered by analyst wixcom ltd nasdaqwix stock holding lessened by legato capital management llc jpmorgan betabuilders japan etf batsbbjp stock holding decreased by madison asset management llc madison asset management llc ha million stock holding in canadian pacific kansa city limited nysecp richardson financial service inc buy new stake in vaneck clo etf nysearcacloi richardson financial service inc cut position in fiserv inc nysefi blackline nasdaqbl stock rating upgraded by da davidson advisor preferred llc take position in the sherwinwilliams company nyseshw legato capital management llc increase stake in spartannash company nasdaqsptn fortis nysefts raised to moderate buy at raymond james financial ralliant nyseral stock rating upgraded by melius research krispy kreme nasdaqdnut upgraded at bnp paribas shake shack shak projected to post quarterly earnings on thursday comparing sma solar technology smtgf and it rival intercontinental hotel group plc lonihg receives consensus recommendation of reduce from analyst contrasting international consolidated airline group otcmktsicagy joby aviation nysejoby hercules capital htgc expected to announce earnings on thursday vericel vcel expected to announce quarterly earnings on thursday fubotv nysefubo and fox nasdaqfox head to head contrast leggett platt leg projected to post earnings on thursday teleflex tfx projected to post quarterly earnings on thursday bayerische motoren werke otcmktsbamxf and nio nysenio financial analysis patrick industry patk projected to post quarterly earnings on thursday dcf claim protocol followed in case of ct boy allegedly held captive gop condemns wait for answer ukraine zelenskyy promise safeguard after street protest over a new anticorruption law ukraine zelenskyy promise safeguard after street protest over a new anticorruption law ukraine zelenskyy promise safeguard after street protest over a new anticorruption law sylogist tsesyz trading up whats next heatwurx nasdaqpcsa share down here what happened proshares ultra vix shortterm future etf batsuvxy trading down here why fundamental global nasdaqfgf stock price down should you sell ct eventim ag co kgaa otcmktscevmf trading higher whats next bosideng international otcmktsbsdgy share down time to sell health science acquisition co nasdaqhsaq stock price up time to buy turk hava yollari anonim ortakligi otcmktstkhvy stock price up should you buy letter another blow dealt against common sense share in hilton grand vacation inc nysehgv purchased by readystate asset management lp share in docusign inc nasdaqdocu bought by efficient advisor llc oak thistle llc ha million stock holding in motorola solution inc nysemsi oak thistle llc buy new position in dynatrace inc nysedt efficient advisor llc buy share of vanguard information technology etf nysearcavgt readystate asset management lp cut stock holding in blackrock capital allocation term trust nysebcat oak thistle llc acquires share of texas roadhouse inc nasdaqtxrh oak thistle llc take position in state street corporation nysestt share in packaging corporation of america nysepkg bought by oak thistle llc oak thistle llc invests million in simpson manufacturing company inc nysessd sym financial corp purchase share of marathon petroleum corporation nysempc oak thistle llc purchase new share in slm corporation nasdaqslm oak thistle llc grows position in viper energy inc nasdaqvnom readystate asset management lp invests in dime community bancshares inc nasdaqdcom south park inadvertently sends liberal conspiracy theory crashing to earth lawyer call fraud on bellring brand but one analyst is optimistic blackstone griddle are up to off during wayfairs black friday in july sale this major retailer is moving into staten island mall empty forever space nuvei debut tool to boost merchant authorization rate irs considers eliminating nonenglish language tax service oxford gray north america corp oxford gray recently filed breach of contract lawsuit for over million against austin texasbased company fintiv also known a mozido owned and led by michael liberty liberty cnbc daily open tesla bad thing come in three brompton fund declares increased etf distribution high dividend bdcs with yield over deltec home launch tranquility a new round hurricaneready home designed for resilience and inspired living bankrupt mexican restaurant chain reopens location after mass closure deltec home launch tranquility a new round hurricaneready home designed for resilience and inspired living bankrupt mexican restaurant chain reopens location after mass closure brookwood investment group llc buy new share in bank of america corporation nysebac brookwood investment group llc invests million in bank of america corporation nysebac nvidia ceo say trump ai plan will fundamentally change u position in year to come mike johnson perfectly explains the flaw in the gop ghislaine maxwell gambit how hulk hogan helped transform pro wrestling from a sideshow to the main event blue ribbon income fund announces monthly distribution federal reserve give journalist a rare tour of building renovation ahead of trump visit federal reserve give journalist a rare tour of building renovation ahead of trump visit unitedhealth say it is under a federal investigation and cooperating unitedhealth say it is under a federal investigation and cooperating unitedhealth say it is under a federal investigation and cooperating unitedhealth say it is under a federal investigation and cooperating unitedhealth say it is under a federal investigation and cooperating unitedhealth say it is under a federal investigation and cooperating unitedhealth say it is under a federal investigation and cooperating unitedhealth say it is under a federal investigation and cooperating unitedhealth say it is under a federal investigation and cooperating unitedhealth say it is under a federal investigation and cooperating unitedhealth say it is under a federal investigation and cooperating unitedhealth say it is under a federal investigation and cooperating unitedhealth say it is under a federal investigation and cooperating blue ribbon income fund announces monthly distribution unitedhealth say it is under a federal investigation and cooperating brompton wellington square aaa clo etf declares distribution wesbanco inc expands presence in northern virginia with commercial loan production office invent inventor develops newly designed head covering plb brompton wellington square aaa clo etf declares distribution construction advance at future johnson city postal center featured to try on your new perplexity pro subscription lake babine nation and west fraser partnership strengthens bcs forest sector inventhelp inventor develops new accessory for roulette player lgt union pacific and norfolk southern confirm merger talk to create coast to coast railroad shenandoah county winery feel squeeze of tariff amplifys ongoing commitment to the pennsylvania association for justice a a justice business partner wesbanco inc expands presence in northern virginia with commercial loan production office construction advance at future johnson city postal center onalaska applebees pair with ihop steve cahalan this health care stock forming a big base breakout is ready to separate from the pack chart indicate why the u fertility rate ha hit an alltime low ibd live qa recap and key stock list for july tesla tumble and alphabet rise to keep wall street near it record tesla tumble and alphabet rise to keep wall street near it record tesla tumble and alphabet rise to keep wall street near it record tesla tumble and alphabet rise to keep wall street near it record feel trapped in a subscription here are your legal right morling financial advisor llc ha million holding in spdr portfolio sp etf nysearcasplg vanguard total stock market etf nysearcavti stock position boosted by morling financial advisor llc america going to win ai race against china trump declares a meta executive applauds strategy how the fed is managing payment new gps kansa city radio station reeling from federal funding cut slash and burn microsoft ceo satya nadella finally address the elephant in the room kind of why is the fed renovation so expensive staff cite tariff and inflation demand optimism lift sugar price pegasystems affirms outlook a cloud strategy power momentum call it the sydney sweeney effect american eagle latest campaign take flight pilgrim to build food plant in northwest georgia ai infrastructure engineer cisco journey strategic wealth llc raise stake in eli lilly and company nyselly eli lilly and company nyselly share bought by sfe investment counsel atrium investment inc buy share of tesla inc nasdaqtsla rtx nysertx post quarterly earnings result beat estimate by eps bg investment service inc acquires share of tesla inc nasdaqtsla ishares msci usa quality factor etf batsqual share sold by benningfield financial advisor llc ninepoint partner lp ha stock holding in eli lilly and company nyselly rtx nysertx update fy earnings guidance northern financial advisor inc increase stake in rtx corporation nysertx northern financial advisor inc sell share of lockheed martin corporation nyselmt birchcreek wealth management llc decrease stock position in invesco qqq nasdaqqqq energy stock worth watching july nd principle wealth partner llc sell share of caterpillar inc nysecat journey strategic wealth llc raise holding in tesla inc nasdaqtsla promising solar stock to follow today july nd diamond hill capital management inc acquires new share in salesforce inc nysecrm share in eli lilly and company nyselly purchased by bg investment service inc invesco qqq nasdaqqqq share sold by journey strategic wealth llc northern financial advisor inc acquires share of rtx corporation nysertx best electric vehicle stock to consider july nd collier list former rite aid location poppi settle class action lawsuit for million poppi settle class action lawsuit for million poppi settle class action lawsuit for million poppi settle class action lawsuit for million cocoa price finish sharply lower a chocolate demand wane unemployment claim in indiana increased last week why your childhood snack are back on shelf cortez community now accepts cryptocurrency from home buyer the latest u cut short gaza ceasefire talk blaming hamas the latest u cut short gaza ceasefire talk blaming hamas the latest u cut short gaza ceasefire talk blaming hamas chicharito face backlash in mexico for very sexist remark on gender role porn racist video force nj election meeting to end a hacker take over gov tony evers announces we wont seek a third term live now texas house hold hearing on congressional redistricting trouble buying it expert warns trump close to triggering second schism with maga base street bordering civic center park to close for red rally activity trump visiting federal reserve amid pressure on powell to cut interest rate humboldt county making progress on ada improvement trump to tour fed ramping up war on central bank trump to tour fed ramping up war on central bank trump to tour fed ramping up war on central bank trump to tour fed ramping up war on central bank trump to tour fed ramping up war on central bank trump to tour fed ramping up war on central bank trump to tour fed ramping up war on central bank trump to tour fed ramping up war on central bank trump to tour fed ramping up war on central bank trump to tour fed ramping up war on central bank trump to tour fed ramping up war on central bank mike johnson perfectly explains the flaw in the gop ghislaine maxwell gambit poll most american approve of recent u supreme court ruling poll most american approve of recent u supreme court ruling poll most american approve of recent u supreme court ruling poll most american approve of recent u supreme court ruling poll most american approve of recent u supreme court ruling poll most american approve of recent u supreme court ruling poll most american approve of recent u supreme court ruling poll most american approve of recent u supreme court ruling poll most american approve of recent u supreme court ruling poll most american approve of recent u supreme court ruling poll most american approve of recent u supreme court ruling poll most american approve of recent u supreme court ruling poll most american approve of recent u supreme court ruling poll most american approve of recent u supreme court ruling poll most american approve of recent u supreme court ruling poll most american approve of recent u supreme court ruling poll most american approve of recent u supreme court ruling poll most american approve of recent u supreme court ruling poll most american approve of recent u supreme court ruling poll most american approve of recent u supreme court ruling poll most american approve of recent u supreme court ruling poll most american approve of recent u supreme court ruling poll most american approve of recent u supreme court ruling poll most american approve of recent u supreme court ruling poll most american approve of recent u supreme court ruling poll most american approve of recent u supreme court ruling poll most american approve of recent u supreme court ruling poll most american approve of recent u supreme court ruling ap news summary at pm edt ap news summary at pm edt u cut short it gaza ceasefire talk and accuses hamas of lacking good faith ghislaine maxwell leaf maximum security florida jail ahead of productive doj meeting federal reserve give journalist a rare tour of building renovation ahead of trump visit eye wide shut hulk hogan became one of trump biggest supporter at rnc after dumping the democrat over a petty slight bar nunn city council hears opposition support for radiant project million of dollar in education funding for minnesota school remain in limbo ahead of school year million of dollar in education funding for minnesota school remain in limbo ahead of school year million of dollar in education funding for minnesota school remain in limbo ahead of school year trump to sign order pushing city and state to remove homeless people from street a list of sport organization that have begun transgender competition ban honor the miracle on ice u olympic team with congressional gold medal a list of sport organization that have begun transgender competition ban zelenskyy move to restore independence of antigraft agency after public anger and eu criticism exbiden chief of staff ron klain cooperative with house gop probe of prezs mental decline national night out returning oct neighborhood invited to register video showing migrant worker moved by forklift prompt action from south korea president pathetic dem rep charged with interfering with ice resort to the only card she ha left top gop lawmaker demand soulsearching from dems amid antiice riot targeting law enforcement alabamian ask prison committee for change in living condition safety standard around the region july best electric vehicle stock to consider july nd oped we support plan for oneway street bus service bike lane and infrastructure why havent epstein file been released most american agree on answer poll find why havent epstein file been released most american agree on answer poll find why havent epstein file been released most american agree on answer poll find why havent epstein file been released most american agree on answer poll find why havent epstein file been released most american agree on answer poll find why havent epstein file been released most american agree on answer poll find ap news summary at pm edt ap news summary at pm edt u cut short it gaza ceasefire talk and accuses hamas of lacking good faith u cut short it gaza ceasefire talk and accuses hamas of lacking good faith u cut short it gaza ceasefire talk and accuses hamas of lacking good faith u cut short it gaza ceasefire talk and accuses hamas of lacking good faith california law requiring background check for ammo is unconstitutional appeal court rule u cut short it gaza ceasefire talk and accuses hamas of lacking good faith u cut short it gaza ceasefire talk and accuses hamas of lacking good faith why havent epstein file been released most american agree on answer poll find declassified intel bombshell cia officer fought brennan on russia collusion doc but he overruled them hindsight is ethic commission investigates oregon housing adviser the dangerous illusion of friendship with qatar and other hostile ally wisconsin dem gov tony evers wont seek third term prediction market forecast socialist landslide for nyc mayor race trump mourns death of hulk hogan he wa maga all the way trump mourns death of hulk hogan he wa maga all the way u recall negotiating team a prospect for a gaza ceasefire dim and humanitarian situation get worse mid america steel employee dy after workplace incident'

$ws.Range("A23").Value = 45863
$ws.Range("A23").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B23").Value = 6388.64013671875
$ws.Range("C23").Value = 6395.81982421875
$ws.Range("D23").Value = 6368.52978515625
$ws.Range("E23").Value = 6370.009765625
$ws.Range("F23").Value = 4470720000
$ws.Range("G23").Value = 0.0039743277792958
$ws.Range("H23").Value = 1
$ws.Range("I23").Value = $headline
$ws.Range("J23").Value = $headline
$ws.Range("K23").Value = 0.16
$ws.Range("L23").Value = 0.064
$ws.Range("M23").Value = 0.776
$ws.Range("N23").Value = 0.9999
$ws.Range("O23").Value = 81
$ws.Range("P23").Value = 70
$ws.Range("Q23").Value = 0
$ws.Range("R23").Value = 397
$ws.Range("S23").Value = 62
$ws.Range("T23").Value = 514
$ws.Range("U23").Value = 110
$ws.Range("V23").Value = 52
$ws.Range("W23").Value = 29
$ws.Range("X23").Value = 201
